$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Bump IG version and regeneration timestamp in the metadata table.
$ws.Range("B3").Value = "0.2.0"
$ws.Range("B8").Value = "2023-10-20T08:59:58+00:00"

# Insert a new "Jurisdiction" property row right after "Contact" (row 10),
# pushing Description/Purpose/Copyright/... down by one row.
$ws.Rows("11:11").Insert()

# The freshly inserted row doesn't carry the data-row look (border/alignment)
# that every other property row has, so clone that formatting from the row
# right below it (now "Description", which already uses the right style)
# before writing the new content into it.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
